$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "Klay Thompson"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Dallas Mavericks"

# Row 4
$ws.Range("A4").Value = "Jaylen Brown"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("C4").Value = "Boston Celtics"

# Row 5
$ws.Range("A5").Value = "Jalen Suggs"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Orlando Magic"

# Row 11
$ws.Range("A11").Value = "Jakob Poeltl"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Toronto Raptors"

# Row 13
$ws.Range("A13").Value = "Clint Capela"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Atlanta Hawks"

# Row 14
$ws.Range("A14").Value = "Chris Paul"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "San Antonio Spurs"

# Row 15
$ws.Range("A15").Value = "Rudy Gobert"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "Minnesota Timberwolves"

# Row 16
$ws.Range("A16").Value = "Jalen Green"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Houston Rockets"
